$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '40.741.08'
$ws.Range("E2").Value = '  +3.40%  '

# Row 3
$ws.Range("D3").Value = '2.214.23'
$ws.Range("E3").Value = '  +2.22%  '

# Row 4
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '229.24'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.13%  '

# Row 6
$ws.Range("E6").Value = '  +1.78%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '64.09'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.40%  '

# Row 8
$ws.Range("E8").Value = '  +0.03%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.406'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.57%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0869'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.70%  '

# Row 11
$ws.Range("E11").Value = '  +0.06%  '

# Row 12
$ws.Range("D12").Value = '2.543.27'
$ws.Range("E12").Value = '  +2.08%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '15.86'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.51%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.33'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.40%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.824'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.04%  '

# Row 16
$ws.Range("E16").Value = '  +0.72%  '

# Row 17
$ws.Range("D17").Value = '2.217.73'
$ws.Range("E17").Value = '  +2.07%  '

# Row 18
$ws.Range("D18").Value = '40.609.24'
$ws.Range("E18").Value = '  +3.02%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '73.85'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.18%  '

# Row 20
$ws.Range("E20").Value = '  +6.29%  '

# Row 21
$ws.Range("E21").Value = '  -0.72%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '250.62'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +7.96%  '

# Row 24
$ws.Range("E24").Value = '  -0.32%  '

# Row 25
$ws.Range("E25").Value = '  -8.33%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.73'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.72%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '173.19'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.48%  '

# Row 28
$ws.Range("E28").Value = '  +1.48%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '20.39'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.77%  '

# Row 30
$ws.Range("E30").Value = '  +2.69%  '

# Row 31
$ws.Range("E31").Value = '  +1.66%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.123'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.14%  '

# Row 33
$ws.Range("E33").Value = '  +0.57%  '

# Row 34
$ws.Range("B34").Value = 'THORChain'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.14'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.12%  '

# Row 35
$ws.Range("B35").Value = 'InternetComputer(DFINITY)'
$ws.Range("C35").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.77'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.34%  '

# Row 36
$ws.Range("E36").Value = '  +1.77%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.82'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.33%  '

# Row 38
$ws.Range("E38").Value = '  +1.51%  '

# Row 39
$ws.Range("E39").Value = '  +0.08%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.84'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +10.76%  '

# Row 41
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.65'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +8.14%  '

# Row 42
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0232'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.74%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '101.56'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.77%  '

# Row 44
$ws.Range("E44").Value = '  +4.29%  '

# Row 45
$ws.Range("D45").Value = '1.521.00'
$ws.Range("E45").Value = '  -1.30%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '17.29'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.04%  '

# Row 47
$ws.Range("E47").Value = '  +1.27%  '

# Row 48
$ws.Range("E48").Value = '  +0.49%  '

# Row 49
$ws.Range("E49").Value = '  +41.00%  '

# Row 50
$ws.Range("E50").Value = '  -0.36%  '

# Row 51
$ws.Range("B51").Value = 'Celestia'
$ws.Range("C51").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.65'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +11.62%  '
